# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-21 14:13:55
# Fix the ordering of "Recorded By" entries in column G of the session analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "backup@backdoor.com, System, system"
    }
}
